# Update absenteeism data rows 2-11 on the active sheet to match the new dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Colaborador_id, B=Colaborador_nome, C=Departamento,
#          D=Motivo_da_ausência, E=Horas_de_ausência, F=Data_da_ausência, G=Salário
$rows = @(
    @{ Row = 2;  A = 27493; B = "Guilherme Monteiro";     C = "P&D";               D = "Doença";              E = 1; F = 45098; G = 8272.280000000001 },
    @{ Row = 3;  A = 50186; B = "Marina Sales";            C = "Marketing";          D = "Outros";              E = 6; F = 45099; G = 7916.48 },
    @{ Row = 4;  A = 69749; B = "Srta. Larissa da Cruz";   C = "Vendas";             D = "Consulta médica";     E = 6; F = 45082; G = 11159.14 },
    @{ Row = 5;  A = 23886; B = "Sr. Fernando Martins";    C = "TI";                 D = "Viagem de negócios";  E = 8; F = 45094; G = 4328.51 },
    @{ Row = 6;  A = 35326; B = "Eduarda Novaes";          C = "Financeiro";         D = "Outros";              E = 5; F = 45097; G = 11254.95 },
    @{ Row = 7;  A = 17095; B = "Nina da Luz";             C = "TI";                 D = "Doença";              E = 6; F = 45083; G = 6257.07 },
    @{ Row = 8;  A = 11142; B = "Ana Sophia Carvalho";     C = "Jurídico";           D = "Doença";              E = 7; F = 45088; G = 5675.45 },
    @{ Row = 9;  A = 18804; B = "Antônio Pereira";         C = "Engenharia";         D = "Doença";              E = 4; F = 45096; G = 3819.41 },
    @{ Row = 10; A = 23592; B = "Luigi Santos";            C = "Recursos Humanos";   D = "Problemas pessoais";  E = 2; F = 45088; G = 11348.25 },
    @{ Row = 11; A = 71361; B = "Luiz Miguel Nunes";       C = "Recursos Humanos";   D = "Viagem de negócios";  E = 2; F = 45104; G = 7685.83 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}
